# Automatische test-sync: 2025-06-26 23:56:50
# Appends the new "Testmail #17" row (row 49) to the Logs sheet, extends the
# conditional-formatting ranges to cover it, and bumps the "Bestelling /
# Levering" tally on the Dashboard sheet from 20 to 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- New row of log data -------------------------------------------------
$ws.Range("A49").Value = "Kun je deze bestelling vandaag verwerken?"
$ws.Range("B49").Value = "mailmind.test@zohomail.eu"
$ws.Range("C49").Value = "Testmail #17: Kun je deze bestelling vandaag verwerken?"
$ws.Range("D49").Value = "Bestelling / Levering"
$ws.Range("E49").Value = "Beste klant,`nHartelijk dank voor uw e-mail. Om de bestelling vandaag te kunnen verwerken, hebben we een ordernummer of klantgegevens nodig. Zou u ons deze gegevens kunnen verstrekken, zodat we uw verzoek verder kunnen onderzoeken?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("F49").Value = "2025-06-26 23:56:10"
$ws.Range("G49").Value = "Ja"
$ws.Range("H49").Value = "Nee"
$ws.Range("I49").Value = "Ja"

# Row 49 holds a multi-line cell (E49); the COM layer auto-expands the row
# height to fit the wrapped text. The source workbook never carries explicit
# row heights, so re-run AutoFit to drop the custom-height flag again and
# keep row 49 consistent with every other row.
$ws.Rows.Item(49).AutoFit()

# --- Extend conditional formatting ranges to include row 49 --------------
function Extend-CF($col) {
    $fcs = $ws.Range($col + "2:" + $col + "48").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($col + "2:" + $col + "49"))
    }
}

Extend-CF("D")
Extend-CF("G")
Extend-CF("H")
Extend-CF("I")

# --- Dashboard: "Bestelling / Levering" count goes from 20 to 21 ---------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 21
